$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells (Coin, Link, Price columns) stay as text even if the
# new value looks like a number (e.g. "1.00", "5.52").
$textCells = @('D2', 'D3', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D32', 'D33', 'D34', 'D35', 'B36', 'C36', 'D36', 'B37', 'C37', 'D37', 'B38', 'C38', 'D38', 'D39', 'B40', 'C40', 'D40', 'B41', 'C41', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'B47', 'C47', 'D47', 'D48', 'B49', 'C49', 'D49', 'D50', 'D51')
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '55.518.04'
$ws.Range('E2').Value = '  +5.21%  '
$ws.Range('D3').Value = '2.490.38'
$ws.Range('E3').Value = '  +6.82%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '486.17'
$ws.Range('E5').Value = '  +10.01%  '
$ws.Range('D6').Value = '139.37'
$ws.Range('E6').Value = '  +14.24%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.36%  '
$ws.Range('D8').Value = '0.509'
$ws.Range('E8').Value = '  +9.10%  '
$ws.Range('D9').Value = '2.488.53'
$ws.Range('E9').Value = '  +7.20%  '
$ws.Range('D10').Value = '0.0979'
$ws.Range('E10').Value = '  +11.58%  '
$ws.Range('D11').Value = '5.52'
$ws.Range('E11').Value = '  +4.80%  '
$ws.Range('D12').Value = '0.330'
$ws.Range('E12').Value = '  +8.97%  '
$ws.Range('E13').Value = '  +2.24%  '
$ws.Range('D14').Value = '2.929.40'
$ws.Range('E14').Value = '  +7.70%  '
$ws.Range('D15').Value = '55.578.58'
$ws.Range('E15').Value = '  +5.38%  '
$ws.Range('D16').Value = '20.72'
$ws.Range('E16').Value = '  +9.60%  '
$ws.Range('D17').Value = '0.0000137'
$ws.Range('E17').Value = '  +14.94%  '
$ws.Range('D18').Value = '2.497.03'
$ws.Range('E18').Value = '  +6.69%  '
$ws.Range('D19').Value = '4.39'
$ws.Range('E19').Value = '  +11.99%  '
$ws.Range('D20').Value = '318.85'
$ws.Range('E20').Value = '  +6.98%  '
$ws.Range('D21').Value = '10.01'
$ws.Range('E21').Value = '  +12.12%  '
$ws.Range('D23').Value = '5.75'
$ws.Range('E23').Value = '  +11.48%  '
$ws.Range('D24').Value = '58.06'
$ws.Range('E24').Value = '  +8.12%  '
$ws.Range('D25').Value = '0.169'
$ws.Range('E25').Value = '  +14.04%  '
$ws.Range('D26').Value = '0.410'
$ws.Range('E26').Value = '  +12.78%  '
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').Value = '2.622.76'
$ws.Range('E28').Value = '  +8.02%  '
$ws.Range('D29').Value = '7.40'
$ws.Range('E29').Value = '  +8.11%  '
$ws.Range('D30').Value = '0.0₃0789'
$ws.Range('E30').Value = '  +16.96%  '
$ws.Range('E31').Value = '  +0.31%  '
$ws.Range('D32').Value = '148.86'
$ws.Range('E32').Value = '  +5.11%  '
$ws.Range('D33').Value = '18.15'
$ws.Range('E33').Value = '  +7.12%  '
$ws.Range('D34').Value = '1.49'
$ws.Range('E34').Value = '  +11.02%  '
$ws.Range('D35').Value = '5.19'
$ws.Range('E35').Value = '  +9.02%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = '0.870'
$ws.Range('E36').Value = '  +4.52%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.12'
$ws.Range('E37').Value = '  +12.79%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '3.68'
$ws.Range('E38').Value = '  +7.14%  '
$ws.Range('D39').Value = '34.15'
$ws.Range('E39').Value = '  +7.15%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').Value = '0.609'
$ws.Range('E41').Value = '  +10.44%  '
$ws.Range('D42').Value = '0.0552'
$ws.Range('E42').Value = '  +9.54%  '
$ws.Range('D43').Value = '3.43'
$ws.Range('E43').Value = '  +9.35%  '
$ws.Range('D44').Value = '1.31'
$ws.Range('E44').Value = '  +8.93%  '
$ws.Range('D45').Value = '4.73'
$ws.Range('E45').Value = '  +10.10%  '
$ws.Range('D46').Value = '1.982.36'
$ws.Range('E46').Value = '  +3.36%  '
$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').Value = '10.12'
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').Value = '0.0908'
$ws.Range('E48').Value = '  +10.27%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').Value = '255.11'
$ws.Range('E49').Value = '  +31.33%  '
$ws.Range('D50').Value = '0.0224'
$ws.Range('E50').Value = '  +9.24%  '
$ws.Range('D51').Value = '17.40'
$ws.Range('E51').Value = '  +11.64%  '
